$d = $word.ActiveDocument

# Locate the "Continuous integration" bullet under "General skills" — note a
# later job-title paragraph reads "Continuous Integration" (capital I), so we
# match case-sensitively to target the right bullet.
$target = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -clike "Continuous integration*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Continuous integration' bullet under General skills"
}

# Insert a new bullet right after it, carrying over the same "Compact" style
# and numbered-list (numId 1002, ilvl 0) formatting.
$target.Range.InsertParagraphAfter() | Out-Null
$newPara = $target.Next()
$newPara.Style = $target.Style
$newPara.Range.ListFormat.ApplyListTemplateWithLevel($target.Range.ListFormat.ListTemplate, $true, 2, $false, 1)
$newPara.Range.Text = "Docker / Rancher container management"
